$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N48").Value = 0.9444444444444444
$ws.Range("M62").Value = 1
$ws.Range("M66").Value = 1
$ws.Range("I69").Value = 0.2380952380952381
$ws.Range("M69").Value = 1
$ws.Range("S69").Value = 0.7619047619047619
$ws.Range("T69").Value = 0.7619047619047619
$ws.Range("M70").Value = 1
$ws.Range("M73").Value = 1
$ws.Range("R73").Value = 0.8432203389830508
$ws.Range("M74").Value = 1
$ws.Range("N74").Value = 1
$ws.Range("O74").Value = 1
$ws.Range("I77").Value = 1
$ws.Range("M77").Value = 1
$ws.Range("S77").Value = 0
$ws.Range("T77").Value = 0
$ws.Range("N78").Value = 1
$ws.Range("O78").Value = 1
$ws.Range("I81").Value = 0.4621212121212122
$ws.Range("M81").Value = 1
$ws.Range("M82").Value = 1
$ws.Range("N82").Value = 0.8888888888888888
$ws.Range("O82").Value = 0.8888888888888888
$ws.Range("I85").Value = 0
$ws.Range("S85").Value = 1
$ws.Range("T85").Value = 1
$ws.Range("M86").Value = 0.9444444444444444
$ws.Range("N86").Value = 0.5555555555555556
$ws.Range("O86").Value = 0.5555555555555556
$ws.Range("M90").Value = 0.8888888888888888
$ws.Range("N90").Value = 0.8888888888888888
$ws.Range("O90").Value = 0.8888888888888888
$ws.Range("M94").Value = 0.9444444444444444
$ws.Range("N94").Value = 0.9444444444444444
$ws.Range("O94").Value = 0.9444444444444444
$ws.Range("M98").Value = 1
$ws.Range("N98").Value = 1
$ws.Range("O98").Value = 1
$ws.Range("M102").Value = 1
$ws.Range("N102").Value = 1
$ws.Range("O102").Value = 1
$ws.Range("M106").Value = 1
$ws.Range("N106").Value = 1
$ws.Range("O106").Value = 1
$ws.Range("M110").Value = 1
$ws.Range("N110").Value = 0.7222222222222222
$ws.Range("O110").Value = 0.7222222222222222
$ws.Range("N114").Value = 0.5555555555555556
$ws.Range("O114").Value = 0.5555555555555556
$ws.Range("I118").Value = 0.4285714285714286
$ws.Range("M118").Value = 1
$ws.Range("N118").Value = 1
$ws.Range("O118").Value = 1
$ws.Range("S118").Value = 1
$ws.Range("T118").Value = 1
$ws.Range("M122").Value = 1
$ws.Range("N122").Value = 0.7222222222222222
$ws.Range("O122").Value = 0.7777777777777778
$ws.Range("S122").Value = 0.425531914893617
$ws.Range("T122").Value = 0.425531914893617
$ws.Range("M126").Value = 1
$ws.Range("S126").Value = 0.5904761904761905
$ws.Range("T126").Value = 0.5904761904761905
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("S127").Value = 1
$ws.Range("T127").Value = 1
$ws.Range("M130").Value = 1
$ws.Range("N130").Value = 0.9444444444444444
$ws.Range("O130").Value = 0.9444444444444444
$ws.Range("S130").Value = 0.6804123711340206
$ws.Range("T130").Value = 0.6804123711340206
$ws.Range("J132").Value = 0
$ws.Range("M134").Value = 1
$ws.Range("N134").Value = 0.8888888888888888
$ws.Range("O134").Value = 0.8333333333333334
$ws.Range("S134").Value = 0.1978021978021978
$ws.Range("T134").Value = 0.1978021978021978
$ws.Range("J136").Value = 0
$ws.Range("I138").Value = 0.3882978723404256
$ws.Range("M138").Value = 1
$ws.Range("R138").Value = 0
$ws.Range("S138").Value = 0.6117021276595744
$ws.Range("T138").Value = 0.6117021276595744
$ws.Range("I142").Value = 0.4545454545454546
$ws.Range("M142").Value = 1
$ws.Range("N142").Value = 1
$ws.Range("O142").Value = 1
$ws.Range("S142").Value = 0.5454545454545454
$ws.Range("T142").Value = 0.5454545454545454
$ws.Range("I144").Value = 1
$ws.Range("S144").Value = 0
$ws.Range("T144").Value = 0
$ws.Range("I145").Value = 1
$ws.Range("J145").Value = 1
$ws.Range("S145").Value = 0
$ws.Range("T145").Value = 0
$ws.Range("N146").Value = 0.9444444444444444
$ws.Range("S146").Value = 0
$ws.Range("T146").Value = 0
$ws.Range("I148").Value = 1
$ws.Range("M148").Value = 1
$ws.Range("S148").Value = 0
$ws.Range("T148").Value = 0
$ws.Range("J149").Value = 1
$ws.Range("M149").Value = 0.9444444444444444
$ws.Range("N150").Value = 0.8888888888888888
$ws.Range("O150").Value = 0.8888888888888888
$ws.Range("R150").Value = 0.3529411764705883
$ws.Range("S150").Value = 0.1176470588235294
$ws.Range("T150").Value = 0.1176470588235294
$ws.Range("I152").Value = 1
$ws.Range("S152").Value = 0
$ws.Range("T152").Value = 0
$ws.Range("M153").Value = 0.9444444444444444
$ws.Range("M154").Value = 1
$ws.Range("N154").Value = 0.8888888888888888
$ws.Range("O154").Value = 0.8888888888888888
$ws.Range("R154").Value = 0.03378378378378377
$ws.Range("I155").Value = 1
$ws.Range("S155").Value = 0
$ws.Range("T155").Value = 0
$ws.Range("I156").Value = 1
$ws.Range("M156").Value = 1
$ws.Range("S156").Value = 0
$ws.Range("T156").Value = 0
$ws.Range("J157").Value = 0
$ws.Range("I158").Value = 0.1111111111111112
$ws.Range("I160").Value = 1
$ws.Range("M160").Value = 1
$ws.Range("N160").Value = 0.8888888888888888
$ws.Range("S160").Value = 0
$ws.Range("T160").Value = 0
$ws.Range("J161").Value = 0
$ws.Range("I162").Value = 0.4252873563218391
$ws.Range("S162").Value = 0.5747126436781609
$ws.Range("T162").Value = 0.5747126436781609
$ws.Range("I164").Value = 1
$ws.Range("S164").Value = 0
$ws.Range("T164").Value = 0
$ws.Range("I165").Value = 1
$ws.Range("S165").Value = 0
$ws.Range("T165").Value = 0
$ws.Range("O166").Value = 1
